# KP-11725 D: Extension of questionnaire's translation files
#
# Adds a new "Variable" column (value "s1" for every data row) right after
# the "Entity Id" column on both the "Translations" sheet and the
# "@@_question" sheet, shifting the existing Type/Index/Original/Translation
# columns one position to the right.

$wb = $excel.ActiveWorkbook

$translations = $wb.Worksheets.Item(1)   # "Translations"
$question     = $wb.Worksheets.Item(2)   # "@@_question"

# ---- Sheet "Translations": insert column B ("Variable" / "s1") ----------
$null = $translations.Columns.Item(2).Insert()

$translations.Range("B1").Value = "Variable"
$translations.Range("B2").Value = "s1"
$translations.Range("B3").Value = "s1"
$translations.Range("B4").Value = "s1"
$translations.Range("B5").Value = "s1"

# ---- Sheet "@@_question": insert column B ("Variable" / "s1") -----------
$null = $question.Columns.Item(2).Insert()

$question.Range("B1").Value = "Variable"
$question.Range("B2").Value = "s1"
# match the text-number-format style already used by the row's Entity Id cell
$question.Range("B2").NumberFormat = "@"

# ---- restore selection / active-sheet state -----------------------------
$null = $question.Select()
$null = $question.Range("B3").Select()

$null = $translations.Select()
$null = $translations.Range("B6").Select()
